# Auto-generated: update market-price derived columns (H-N) for specific rows
# across multiple sheets, per scheduled-runner price refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 7683.3335
$ws.Range("I19").Value = 25712.75
$ws.Range("K19").Value = 25712.75
$ws.Range("M19").Value = -25537.75
$ws.Range("H92").Value = 426.77777
$ws.Range("I92").Value = 418.25
$ws.Range("K92").Value = 418.25
$ws.Range("M92").Value = 829.75
$ws.Range("H129").Value = 867.3226
$ws.Range("I129").Value = 542.5
$ws.Range("J129").Value = 889.7241
$ws.Range("K129").Value = 1627.5
$ws.Range("L129").Value = 2669.1723
$ws.Range("M129").Value = 3372.5
$ws.Range("N129").Value = -12669.1723
$ws.Range("H137").Value = 1926.0769
$ws.Range("I137").Value = 1435
$ws.Range("J137").Value = 2595.7273
$ws.Range("K137").Value = 4305
$ws.Range("L137").Value = 7787.1819
$ws.Range("M137").Value = -1755
$ws.Range("N137").Value = -12887.1819
$ws.Range("H138").Value = 1554.2325
$ws.Range("I138").Value = 1204.0385
$ws.Range("J138").Value = 2089.8235
$ws.Range("K138").Value = 3612.1155
$ws.Range("L138").Value = 6269.470499999999
$ws.Range("M138").Value = 1527.8845
$ws.Range("N138").Value = -16549.4705

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2711.67
$ws.Range("I32").Value = 1984.011
$ws.Range("J32").Value = 10069.111
$ws.Range("K32").Value = 1984.011
$ws.Range("L32").Value = 10069.111
$ws.Range("M32").Value = -1697.011
$ws.Range("N32").Value = -10643.111
$ws.Range("H61").Value = 2573.6191
$ws.Range("I61").Value = 1949.7646
$ws.Range("J61").Value = 5225
$ws.Range("K61").Value = 1949.7646
$ws.Range("L61").Value = 5225
$ws.Range("M61").Value = -1737.7646
$ws.Range("N61").Value = -5649
$ws.Range("H74").Value = 1386.2354
$ws.Range("I74").Value = 1352.75
$ws.Range("J74").Value = 1416
$ws.Range("K74").Value = 1352.75
$ws.Range("L74").Value = 1416
$ws.Range("M74").Value = -478.75
$ws.Range("N74").Value = -3164
$ws.Range("H77").Value = 1386.2354
$ws.Range("I77").Value = 1352.75
$ws.Range("J77").Value = 1416
$ws.Range("K77").Value = 6763.75
$ws.Range("L77").Value = 7080
$ws.Range("M77").Value = -2395.75
$ws.Range("N77").Value = -15816
$ws.Range("H88").Value = 27501.5
$ws.Range("I88").Value = 100006
$ws.Range("J88").Value = 3333.3333
$ws.Range("K88").Value = 100006
$ws.Range("L88").Value = 3333.3333
$ws.Range("M88").Value = -99600
$ws.Range("N88").Value = -4145.3333
$ws.Range("H91").Value = 27501.5
$ws.Range("I91").Value = 100006
$ws.Range("J91").Value = 3333.3333
$ws.Range("K91").Value = 100006
$ws.Range("L91").Value = 3333.3333
$ws.Range("M91").Value = -98602
$ws.Range("N91").Value = -6141.3333
$ws.Range("H132").Value = 3071.76
$ws.Range("I132").Value = 1217.8667
$ws.Range("J132").Value = 5852.6
$ws.Range("K132").Value = 3653.6001
$ws.Range("L132").Value = 17557.8
$ws.Range("M132").Value = -1123.6001
$ws.Range("N132").Value = -22617.8
$ws.Range("H136").Value = 2573.6191
$ws.Range("I136").Value = 1949.7646
$ws.Range("J136").Value = 5225
$ws.Range("K136").Value = 5849.293799999999
$ws.Range("L136").Value = 15675
$ws.Range("M136").Value = -3299.293799999999
$ws.Range("N136").Value = -20775

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2062.0244
$ws.Range("I134").Value = 903.88464
$ws.Range("K134").Value = 2711.65392
$ws.Range("M134").Value = -176.6539199999997

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2229.0476
$ws.Range("I31").Value = 2121.5789
$ws.Range("J31").Value = 3250
$ws.Range("K31").Value = 2121.5789
$ws.Range("L31").Value = 3250
$ws.Range("M31").Value = -1826.5789
$ws.Range("N31").Value = -3840
$ws.Range("H34").Value = 2229.0476
$ws.Range("I34").Value = 2121.5789
$ws.Range("J34").Value = 3250
$ws.Range("K34").Value = 2121.5789
$ws.Range("L34").Value = 3250
$ws.Range("M34").Value = -1919.5789
$ws.Range("N34").Value = -3654
$ws.Range("H58").Value = 1763.325
$ws.Range("I58").Value = 941.4231
$ws.Range("J58").Value = 3289.7144
$ws.Range("K58").Value = 941.4231
$ws.Range("L58").Value = 3289.7144
$ws.Range("M58").Value = -738.4231
$ws.Range("N58").Value = -3695.7144
$ws.Range("H132").Value = 2329.3125
$ws.Range("I132").Value = 1516.6666
$ws.Range("K132").Value = 4549.9998
$ws.Range("M132").Value = -2019.9998
$ws.Range("H134").Value = 3617.8572
$ws.Range("I134").Value = 3335.0908
$ws.Range("J134").Value = 4654.6665
$ws.Range("K134").Value = 10005.2724
$ws.Range("L134").Value = 13963.9995
$ws.Range("M134").Value = -7470.2724
$ws.Range("N134").Value = -19033.9995
$ws.Range("H136").Value = 1763.325
$ws.Range("I136").Value = 941.4231
$ws.Range("J136").Value = 3289.7144
$ws.Range("K136").Value = 2824.2693
$ws.Range("L136").Value = 9869.143199999999
$ws.Range("M136").Value = -274.2692999999999
$ws.Range("N136").Value = -14969.1432

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1154.72
$ws.Range("I131").Value = 788.5
$ws.Range("J131").Value = 1164.7534
$ws.Range("K131").Value = 2365.5
$ws.Range("L131").Value = 3494.2602
$ws.Range("M131").Value = 2674.5
$ws.Range("N131").Value = -13574.2602

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 4018.5
$ws.Range("J132").Value = 4539.385
$ws.Range("L132").Value = 13618.155
$ws.Range("N132").Value = -18678.155

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 24430.182
$ws.Range("I40").Value = 34060.934
$ws.Range("J40").Value = 3792.8572
$ws.Range("K40").Value = 34060.934
$ws.Range("L40").Value = 3792.8572
$ws.Range("M40").Value = -33924.934
$ws.Range("N40").Value = -4064.8572
$ws.Range("H46").Value = 113111.22
$ws.Range("I46").Value = 202800.2
$ws.Range("J46").Value = 1000
$ws.Range("K46").Value = 202800.2
$ws.Range("L46").Value = 1000
$ws.Range("M46").Value = -202612.2
$ws.Range("N46").Value = -1376
$ws.Range("H132").Value = 6565.7256
$ws.Range("I132").Value = 9161.464
$ws.Range("J132").Value = 3405.6956
$ws.Range("K132").Value = 27484.392
$ws.Range("L132").Value = 10217.0868
$ws.Range("M132").Value = -24954.392
$ws.Range("N132").Value = -15277.0868
$ws.Range("H136").Value = 7365.273
$ws.Range("I136").Value = 9681.6
$ws.Range("J136").Value = 5435
$ws.Range("K136").Value = 29044.8
$ws.Range("L136").Value = 16305
$ws.Range("M136").Value = -26494.8
$ws.Range("N136").Value = -21405

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1996.7632
$ws.Range("I132").Value = 1583.875
$ws.Range("J132").Value = 4198.8335
$ws.Range("K132").Value = 4751.625
$ws.Range("L132").Value = 12596.5005
$ws.Range("M132").Value = -2221.625
$ws.Range("N132").Value = -17656.5005
$ws.Range("H136").Value = 5961.3706
$ws.Range("I136").Value = 1114.7
$ws.Range("J136").Value = 8812.352999999999
$ws.Range("K136").Value = 3344.1
$ws.Range("L136").Value = 26437.059
$ws.Range("M136").Value = -794.1000000000004
$ws.Range("N136").Value = -31537.059
